# Applies the "Latest data" refresh to docs/particelle_non_trovate.xlsx.
# A block of 50 new "codice_particella" rows (all belonging to
# codice_comune_catastale 362) is inserted right after the existing row 12
# (i.e. before the current row 13), and every row that used to follow is
# pushed down by 50 rows while keeping its own data intact. The running
# index in column A (which is a plain literal, not a formula) is then
# renumbered so it keeps counting up without gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. The 50 new data rows that get inserted -----------------------------
# Columns: [ index (col A), codice_particella (col B), codice_comune_catastale (col C) ]
$newData = @(
    @(11, "240", 362),
    @(12, "269/1", 362),
    @(13, "269/22", 362),
    @(14, "270/5", 362),
    @(15, "279", 362),
    @(16, "288", 362),
    @(17, "333/1", 362),
    @(18, "333/2", 362),
    @(19, "333/3", 362),
    @(20, "409", 362),
    @(21, "812/1", 362),
    @(22, "856/1", 362),
    @(23, "862", 362),
    @(24, "863", 362),
    @(25, "864/1", 362),
    @(26, "864/2", 362),
    @(27, "1029", 362),
    @(28, "1254", 362),
    @(29, "1256/3", 362),
    @(30, "1256/4", 362),
    @(31, "1305", 362),
    @(32, "1343", 362),
    @(33, "1344", 362),
    @(34, "1345", 362),
    @(35, "1346", 362),
    @(36, "1347/1", 362),
    @(37, "1347/2", 362),
    @(38, "1348", 362),
    @(39, "1383", 362),
    @(40, "1384", 362),
    @(41, "1830/1", 362),
    @(42, "1830/2", 362),
    @(43, "1830/6", 362),
    @(44, "1830/13", 362),
    @(45, "1830/14", 362),
    @(46, "1830/18", 362),
    @(47, "1830/25", 362),
    @(48, "1830/32", 362),
    @(49, "1830/34", 362),
    @(50, "1830/35", 362),
    @(51, "1831", 362),
    @(52, "1832", 362),
    @(53, "1833", 362),
    @(54, "1834", 362),
    @(55, "1836", 362),
    @(56, "1838", 362),
    @(57, ".158", 362),
    @(58, ".425", 362),
    @(59, ".436", 362),
    @(60, ".458", 362)
)

$numNew = $newData.Count
$startRow = 13
$endRow = $wb.Worksheets.Item(1).UsedRange.Rows.Count

# --- 2. Insert blank rows for the new data ---------------------------------
$ws.Rows(($startRow.ToString() + ":" + ($startRow + $numNew - 1).ToString())).Insert()

# Re-apply column A's formatting (bold, centered, thin border) to the freshly
# inserted cells -- Insert() alone does not faithfully reuse the existing
# style for the new row.
$newARange = $ws.Range("A" + $startRow + ":A" + ($startRow + $numNew - 1))
$newARange.Font.Bold = $true
$newARange.HorizontalAlignment = -4108
$newARange.VerticalAlignment = -4160
$newARange.Borders.LineStyle = 1

# Force column B of the new block to be stored as text so that purely
# numeric-looking particella codes (e.g. "240", ".158") are not silently
# turned into numbers, matching how the rest of the sheet stores this column.
$newBRange = $ws.Range("B" + $startRow + ":B" + ($startRow + $numNew - 1))
$newBRange.NumberFormat = "@"

# --- 3. Write the new rows ---------------------------------------------------
for ($i = 0; $i -lt $numNew; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- 4. Renumber column A for all the rows that were pushed down -----------
# They kept their original index value after the Insert(); bump it by the
# number of newly inserted rows so the running count stays continuous.
$shiftedFirst = $startRow + $numNew
for ($r = $shiftedFirst; $r -le $endRow + $numNew; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    if ($old -ne $null -and $old -ne "") {
        $cell.Value = [double]$old + $numNew
    }
}
